$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GC_Codes")

$newCodes = @(
    "LZOKCKSKITLTMZPH ",
    "LAHZTYLVDRPJLSRR ",
    "FLSJSZWQDIYJRCTH ",
    "WOODRMOLGSQHYYGF ",
    "ZYAJYSVDTCRCTWYD ",
    "YOJOKYYVHQLMIJFW "
)

for ($i = 0; $i -lt $newCodes.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $newCodes[$i]
}

$ws.Range("A8:A11").ClearContents()
$ws.Range("A8:A11").EntireRow.Delete()

$ws.Activate()
$ws.Range("A8").Select()
